$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text number-format on every Price (D) cell we touch so Excel's
# COM auto-type-inference doesn't coerce numeric-looking strings (e.g.
# '0.530', '166.10') into floating point Doubles and lose trailing zeros.
$priceCells = @('D2', 'D3', 'D5', 'D6', 'D8', 'D9', 'D10', 'D12', 'D14', 'D15', 'D17', 'D18', 'D19', 'D20', 'D21', 'D25', 'D26', 'D28', 'D30', 'D31', 'D32', 'D33', 'D37', 'D38', 'D39', 'D41', 'D43', 'D46', 'D47', 'D48', 'D49')
foreach ($addr in $priceCells) { $ws.Range($addr).NumberFormat = "@" }

$ws.Range("D2").Value = '70.264.62'
$ws.Range("E2").Value = '  -0.34%  '

$ws.Range("D3").Value = '2.520.87'
$ws.Range("E3").Value = '  -1.24%  '

$ws.Range("E4").Value = '  +0.08%  '

$ws.Range("D5").Value = '575.92'
$ws.Range("E5").Value = '  -0.34%  '

$ws.Range("D6").Value = '166.10'
$ws.Range("E6").Value = '  -2.27%  '

$ws.Range("D8").Value = '0.521'
$ws.Range("E8").Value = '  +2.19%  '

$ws.Range("D9").Value = '2.523.93'
$ws.Range("E9").Value = '  -1.04%  '

$ws.Range("D10").Value = '0.161'
$ws.Range("E10").Value = '  -2.52%  '

$ws.Range("E11").Value = '  -1.02%  '

$ws.Range("D12").Value = '0.356'
$ws.Range("E12").Value = '  +3.07%  '

$ws.Range("E13").Value = '  +1.84%  '

$ws.Range("D14").Value = '2.986.85'
$ws.Range("E14").Value = '  -1.05%  '

$ws.Range("D15").Value = '70.187.78'
$ws.Range("E15").Value = '  -0.25%  '

$ws.Range("E16").Value = '  -2.04%  '

$ws.Range("D17").Value = '25.03'
$ws.Range("E17").Value = '  -0.34%  '

$ws.Range("D18").Value = '2.529.55'
$ws.Range("E18").Value = '  -0.96%  '

$ws.Range("D19").Value = '11.47'
$ws.Range("E19").Value = '  -2.04%  '

$ws.Range("D20").Value = '7.82'
$ws.Range("E20").Value = '  +1.44%  '

$ws.Range("D21").Value = '351.76'
$ws.Range("E21").Value = '  -2.72%  '

$ws.Range("E22").Value = '  -0.99%  '

$ws.Range("E23").Value = '  -1.31%  '

$ws.Range("E24").Value = '  +0.02%  '

$ws.Range("D25").Value = '70.45'
$ws.Range("E25").Value = '  +0.66%  '

$ws.Range("D26").Value = '4.01'
$ws.Range("E26").Value = '  -1.97%  '

$ws.Range("D28").Value = '8.87'
$ws.Range("E28").Value = '  -4.77%  '

$ws.Range("E29").Value = '  +0.74%  '

$ws.Range("D30").Value = '0.0₃0900'
$ws.Range("E30").Value = '  -2.92%  '

$ws.Range("D31").Value = '7.89'
$ws.Range("E31").Value = '  +0.33%  '

$ws.Range("D32").Value = '465.55'
$ws.Range("E32").Value = '  -4.28%  '

$ws.Range("D33").Value = '1.25'
$ws.Range("E33").Value = '  -3.83%  '

$ws.Range("E34").Value = '  -1.42%  '

$ws.Range("E35").Value = '  +0.28%  '

$ws.Range("E36").Value = '  +0.87%  '

$ws.Range("D37").Value = '155.71'
$ws.Range("E37").Value = '  -0.75%  '

$ws.Range("D38").Value = '19.07'
$ws.Range("E38").Value = '  +1.29%  '

$ws.Range("D39").Value = '18.67'
$ws.Range("E39").Value = '  -0.20%  '

$ws.Range("D41").Value = '4.79'
$ws.Range("E41").Value = '  +0.33%  '

$ws.Range("E42").Value = '  -0.73%  '

$ws.Range("D43").Value = '1.61'
$ws.Range("E43").Value = '  -4.25%  '

$ws.Range("E44").Value = '  +0.04%  '

$ws.Range("E45").Value = '  -13.68%  '

$ws.Range("D46").Value = '2.30'
$ws.Range("E46").Value = '  -7.08%  '

$ws.Range("D47").Value = '142.97'
$ws.Range("E47").Value = '  -1.49%  '

$ws.Range("D48").Value = '0.530'
$ws.Range("E48").Value = '  -0.51%  '

$ws.Range("D49").Value = '3.49'
$ws.Range("E49").Value = '  -1.81%  '

$ws.Range("E50").Value = '  -3.21%  '

$ws.Range("E51").Value = '  -0.91%  '

# Restore default (General) styling on the Price cells so the workbook's
# style table isn't left referencing a Text number format on these cells.
foreach ($addr in $priceCells) { $ws.Range($addr).Style = "Normal" }
